$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for rows 2-31, column G
$kValues = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 2
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 1
    30 = 1
    31 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
